$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 53.7
$ws.Range("C2").Value = 23.765

$ws.Range("B3").Value = 55.1
$ws.Range("C3").Value = 24.23

$ws.Range("B4").Value = 52.3
$ws.Range("C4").Value = 24.02

$ws.Range("B5").Value = 54.8
$ws.Range("C5").Value = 24.318

$ws.Range("B6").Value = 57.3
$ws.Range("C6").Value = 23.885

$ws.Range("B7").Value = 56.1
$ws.Range("C7").Value = 23.5

$ws.Range("B8").Value = 56.5
$ws.Range("C8").Value = 23.686

$ws.Range("B9").Value = 54.4
$ws.Range("C9").Value = 23.903
